$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.957
$ws.Range("C2").Value = 0.83
$ws.Range("E2").Value = 0.968
$ws.Range("F2").Value = 0.871
$ws.Range("H2").Value = 0.18
$ws.Range("I2").Value = 0.377
$ws.Range("J2").Value = 0.45
$ws.Range("K2").Value = 0.032
$ws.Range("L2").Value = 0.133
$ws.Range("M2").Value = 0.165
$ws.Range("N2").Value = 0.07699999958276749
$ws.Range("P2").Value = 0.1099999994039536
$ws.Range("Q2").Value = 0.06400000303983688
$ws.Range("R2").Value = 0.07699999958276749
$ws.Range("S2").Value = 0.0820000022649765
$ws.Range("B3").Value = 0.567
$ws.Range("C3").Value = 0.503
$ws.Range("D3").Value = 0.589
$ws.Range("E3").Value = 0.673
$ws.Range("F3").Value = 0.622
$ws.Range("G3").Value = 0.688
$ws.Range("H3").Value = 0.664
$ws.Range("I3").Value = 0.745
$ws.Range("J3").Value = 0.687
$ws.Range("K3").Value = 0.361
$ws.Range("L3").Value = 0.431
$ws.Range("M3").Value = 0.36
$ws.Range("N3").Value = 0.1630000025033951
$ws.Range("O3").Value = 0.1759999990463257
$ws.Range("P3").Value = 0.1609999984502792
$ws.Range("Q3").Value = 0.1180000007152557
$ws.Range("R3").Value = 0.1280000060796738
$ws.Range("S3").Value = 0.1150000020861626
$ws.Range("B4").Value = 0.373
$ws.Range("C4").Value = 0.37
$ws.Range("D4").Value = 0.409
$ws.Range("E4").Value = 0.494
$ws.Range("F4").Value = 0.48
$ws.Range("G4").Value = 0.52
$ws.Range("H4").Value = 0.896
$ws.Range("I4").Value = 0.894
$ws.Range("J4").Value = 0.919
$ws.Range("K4").Value = 0.602
$ws.Range("L4").Value = 0.613
$ws.Range("M4").Value = 0.6
$ws.Range("N4").Value = 0.2189999967813492
$ws.Range("O4").Value = 0.2179999947547913
$ws.Range("Q4").Value = 0.1720000058412552
$ws.Range("R4").Value = 0.1739999949932098
$ws.Range("S4").Value = 0.1710000038146973
$ws.Range("B5").Value = 0.419
$ws.Range("C5").Value = 0.446
$ws.Range("D5").Value = 0.375
$ws.Range("E5").Value = 0.561
$ws.Range("F5").Value = 0.578
$ws.Range("H5").Value = 0.91
$ws.Range("I5").Value = 0.94
$ws.Range("J5").Value = 0.981
$ws.Range("K5").Value = 0.556
$ws.Range("L5").Value = 0.556
$ws.Range("M5").Value = 0.621
$ws.Range("N5").Value = 0.2280000001192093
$ws.Range("O5").Value = 0.2240000069141388
$ws.Range("P5").Value = 0.2430000007152557
$ws.Range("Q5").Value = 0.1630000025033951
$ws.Range("R5").Value = 0.1560000032186508
$ws.Range("S5").Value = 0.1770000010728836
$ws.Range("B6").Value = 0.9
$ws.Range("C6").Value = 0.713
$ws.Range("D6").Value = 0.717
$ws.Range("E6").Value = 0.925
$ws.Range("F6").Value = 0.783
$ws.Range("G6").Value = 0.788
$ws.Range("H6").Value = 0.285
$ws.Range("I6").Value = 0.492
$ws.Range("J6").Value = 0.533
$ws.Range("K6").Value = 0.077
$ws.Range("L6").Value = 0.225
$ws.Range("M6").Value = 0.233
$ws.Range("N6").Value = 0.08900000154972076
$ws.Range("O6").Value = 0.1270000040531158
$ws.Range("P6").Value = 0.1289999932050705
$ws.Range("Q6").Value = 0.07000000029802322
$ws.Range("R6").Value = 0.09200000017881393
$ws.Range("S6").Value = 0.09300000220537186
$ws.Range("B7").Value = 0.59
$ws.Range("C7").Value = 0.548
$ws.Range("D7").Value = 0.541
$ws.Range("E7").Value = 0.692
$ws.Range("F7").Value = 0.658
$ws.Range("G7").Value = 0.656
$ws.Range("H7").Value = 0.649
$ws.Range("I7").Value = 0.665
$ws.Range("J7").Value = 0.693
$ws.Range("K7").Value = 0.346
$ws.Range("L7").Value = 0.375
$ws.Range("M7").Value = 0.387
$ws.Range("N7").Value = 0.1599999964237213
$ws.Range("O7").Value = 0.1669999957084656
$ws.Range("P7").Value = 0.1700000017881393
$ws.Range("Q7").Value = 0.1159999966621399
$ws.Range("R7").Value = 0.1260000020265579
$ws.Range("S7").Value = 0.125
$ws.Range("A8").Value = "40_fnn"
$ws.Range("B8").Value = 0.36
$ws.Range("C8").Value = 0.339
$ws.Range("D8").Value = 0.375
$ws.Range("E8").Value = 0.51
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.523
$ws.Range("H8").Value = 0.926
$ws.Range("I8").Value = 0.922
$ws.Range("J8").Value = 0.904
$ws.Range("K8").Value = 0.612
$ws.Range("L8").Value = 0.617
$ws.Range("M8").Value = 0.589
$ws.Range("N8").Value = 0.2249999940395355
$ws.Range("O8").Value = 0.2310000061988831
$ws.Range("P8").Value = 0.2240000069141388
$ws.Range("Q8").Value = 0.1729999929666519
$ws.Range("R8").Value = 0.1749999970197678
$ws.Range("S8").Value = 0.1710000038146973
$ws.Range("U8").Value = "fnn"
$ws.Range("A9").Value = "40_gat"
$ws.Range("B9").Value = 0.467
$ws.Range("C9").Value = 0.356
$ws.Range("D9").Value = 0.37
$ws.Range("E9").Value = 0.596
$ws.Range("F9").Value = 0.508
$ws.Range("G9").Value = 0.519
$ws.Range("H9").Value = 0.885
$ws.Range("I9").Value = 1.194
$ws.Range("J9").Value = 0.995
$ws.Range("K9").Value = 0.525
$ws.Range("L9").Value = 0.758
$ws.Range("M9").Value = 0.637
$ws.Range("N9").Value = 0.2160000056028366
$ws.Range("O9").Value = 0.300000011920929
$ws.Range("P9").Value = 0.2419999986886978
$ws.Range("Q9").Value = 0.1550000011920929
$ws.Range("R9").Value = 0.2080000042915344
$ws.Range("S9").Value = 0.1780000030994415
$ws.Range("U9").Value = "gat"
$ws.Range("B10").Value = 0.659
$ws.Range("C10").Value = 0.676
$ws.Range("D10").Value = 0.599
$ws.Range("E10").Value = 0.742
$ws.Range("F10").Value = 0.767
$ws.Range("G10").Value = 0.696
$ws.Range("H10").Value = 0.725
$ws.Range("I10").Value = 0.483
$ws.Range("J10").Value = 0.706
$ws.Range("K10").Value = 0.342
$ws.Range("L10").Value = 0.233
$ws.Range("M10").Value = 0.363
$ws.Range("N10").Value = 0.1689999997615814
$ws.Range("O10").Value = 0.1150000020861626
$ws.Range("P10").Value = 0.1700000017881393
$ws.Range("Q10").Value = 0.1129999980330467
$ws.Range("R10").Value = 0.08399999886751175
$ws.Range("S10").Value = 0.1190000027418137
$ws.Range("B11").Value = 0.627
$ws.Range("C11").Value = 0.446
$ws.Range("D11").Value = 0.494
$ws.Range("E11").Value = 0.717
$ws.Range("F11").Value = 0.6
$ws.Range("G11").Value = 0.619
$ws.Range("H11").Value = 0.677
$ws.Range("I11").Value = 0.632
$ws.Range("J11").Value = 0.719
$ws.Range("K11").Value = 0.342
$ws.Range("L11").Value = 0.4
$ws.Range("M11").Value = 0.425
$ws.Range("N11").Value = 0.1580000072717667
$ws.Range("O11").Value = 0.1529999971389771
$ws.Range("P11").Value = 0.1749999970197678
$ws.Range("Q11").Value = 0.1140000000596046
$ws.Range("R11").Value = 0.1220000013709068
$ws.Range("S11").Value = 0.1309999972581863
$ws.Range("A12").Value = "10_gat"
$ws.Range("B12").Value = 0.475
$ws.Range("C12").Value = 0.451
$ws.Range("D12").Value = 0.38
$ws.Range("E12").Value = 0.592
$ws.Range("F12").Value = 0.6
$ws.Range("G12").Value = 0.521
$ws.Range("H12").Value = 1.033
$ws.Range("I12").Value = 0.775
$ws.Range("J12").Value = 1.052
$ws.Range("K12").Value = 0.6
$ws.Range("L12").Value = 0.467
$ws.Range("M12").Value = 0.667
$ws.Range("N12").Value = 0.2509999871253967
$ws.Range("O12").Value = 0.1959999948740005
$ws.Range("P12").Value = 0.257999986410141
$ws.Range("Q12").Value = 0.1700000017881393
$ws.Range("R12").Value = 0.1340000033378601
$ws.Range("S12").Value = 0.1860000044107437
$ws.Range("U12").Value = "gat"
$ws.Range("A13").Value = "10_fnn"
$ws.Range("B13").Value = 0.325
$ws.Range("C13").Value = 0.368
$ws.Range("D13").Value = 0.355
$ws.Range("E13").Value = 0.425
$ws.Range("F13").Value = 0.467
$ws.Range("G13").Value = 0.435
$ws.Range("H13").Value = 1.372
$ws.Range("I13").Value = 1.354
$ws.Range("J13").Value = 1.317
$ws.Range("K13").Value = 0.967
$ws.Range("L13").Value = 0.9
$ws.Range("M13").Value = 0.912
$ws.Range("N13").Value = 0.3339999914169312
$ws.Range("O13").Value = 0.3269999921321869
$ws.Range("P13").Value = 0.3210000097751617
$ws.Range("Q13").Value = 0.2569999992847443
$ws.Range("R13").Value = 0.2339999973773956
$ws.Range("S13").Value = 0.2450000047683716
$ws.Range("U13").Value = "fnn"
